$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the paired "_old" / "_new" column-header suffixes to the
#    format-version-specific suffixes "_FV2210" / "_FV2304".
#    Columns A-J carry the "_old" (-> "_FV2210") headers, column K is the
#    untouched "diff" header, and columns L-U carry the "_new"
#    (-> "_FV2304") headers.
# ---------------------------------------------------------------------------
$oldHeaders = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$leftCols  = @("A","B","C","D","E","F","G","H","I","J")
$rightCols = @("L","M","N","O","P","Q","R","S","T","U")

for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $ws.Range($leftCols[$i]  + "1").Value = $oldHeaders[$i] + "_FV2210"
    $ws.Range($rightCols[$i] + "1").Value = $oldHeaders[$i] + "_FV2304"
}

# ---------------------------------------------------------------------------
# 2) Turn the used range into a native Excel table ("Table1") without
#    disturbing the existing header-row formatting (bold font, grey fill,
#    thin border, centered + wrapped alignment) and without pulling in a
#    table-style header override. To do this, the current header formatting
#    is stashed in a scratch row, the header format is reset to the
#    workbook default before the table is created (so Excel does not bake a
#    header dxf into the table definition), and the original formatting is
#    restored onto the header afterwards.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A200:U200")

$headerRange.Copy()
$scratchRange.PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange.Style = "Normal"

$tableRange = $ws.Range("A1:U90")
$tbl = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$tbl.TableStyle = ""

$scratchRange.Copy()
$headerRange.PasteSpecial(-4122)    # xlPasteFormats
$excel.CutCopyMode = $false

$scratchRange.Clear()

# ---------------------------------------------------------------------------
# 3) Freeze the header row.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "done"
